$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.273.62'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').Value = '2.556.29'
$ws.Range('E3').Value = '  -3.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '515.12'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.87'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.560'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').Value = '2.568.08'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.58'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.83%  '
$ws.Range('E11').Value = '  -2.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.323'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -4.23%  '
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').Value = '3.011.89'
$ws.Range('E14').Value = '  -3.57%  '
$ws.Range('D15').Value = '57.287.25'
$ws.Range('E15').Value = '  -2.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.16'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.10%  '
$ws.Range('E17').Value = '  -2.62%  '
$ws.Range('D18').Value = '2.526.14'
$ws.Range('E18').Value = '  -5.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '337.48'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('E20').Value = '  -2.63%  '
$ws.Range('E21').Value = '  -2.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.24'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.47'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.400'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.94%  '
$ws.Range('D28').Value = '2.678.82'
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.93'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.36%  '
$ws.Range('E30').Value = '  -6.81%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.40'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.02%  '
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.58'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '149.33'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.99'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.16%  '
$ws.Range('E37').Value = '  -4.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.863'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.03'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('E40').Value = '  -5.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.44'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.52'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '269.19'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.97%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.65'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0952'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.585'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.75'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0520'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.71%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.971.08'
$ws.Range('E50').Value = '  -4.09%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.36'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.74%  '
